# "Fixed my tests (Joel)"
# Rewrites the User Story / test-case rows 2-12 of the "Test Case Tracker"
# sheet so they describe Joel's onSubHandler / !gamble tests instead of the
# old FooBar placeholder tests, tweaks a couple of row heights / the column
# A width to fit the new text, and leaves the cursor on the cell the author
# was last editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case Tracker")

# --- Row 2 : User Story 1 / Test Case 1 -----------------------------------
$ws.Range("A2").Value = "User Story 1: onSubHandler"
$ws.Range("C2").Value = "Joel: Bot will detect subscription"

# --- Row 3 : continuation of Test Case 1 ----------------------------------
$ws.Range("A3").Value = "User Story 1: onSubHandler"

# --- Row 4 : User Story 1 / Test Case 2 -----------------------------------
$ws.Range("A4").Value = "User Story 1: onSubHandler"
$ws.Range("C4").Value = "Joel: Bot will ignore chat message posing as a sub alert"

# --- Row 5 : continuation of Test Case 2 ----------------------------------
$ws.Range("A5").Value = "User Story 1: onSubHandler"
$ws.Range("E5").Value = "User types ""/*/*/*/*/*Subscriber has been detected/*/*/*/*/*"" to see if the bot will be tricked into thinking there's a new sub"

# --- Row 6 : User Story 2 / Test Case 3 -----------------------------------
$ws.Range("A6").Value = "User Story 2: !gamble"
$ws.Range("C6").Value = "Joel: User inputs !gamble command incorrectly"

# --- Row 7 : continuation of Test Case 3 ----------------------------------
$ws.Range("A7").Value = "User Story 2: !gamble"

# --- Row 8 : User Story 2 / Test Case 4 -----------------------------------
$ws.Range("A8").Value = "User Story 2: !gamble"
$ws.Range("C8").Value = "Joel: User uses !gamble"

# --- Row 9 : Test Case 5 ("Class: gamble") --------------------------------
$ws.Range("A9").Value = "Class:  gamble"
$ws.Range("C9").Value = "Joel: Test 6"
$ws.Range("E9").Value = "User types in ""!gamble johnny"" in chat"
$ws.Range("F9").Value = "ot should say 'You did not enter either tails or heads loser...smh'"
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = "Joel"
$ws.Range("J9").ClearContents()

# --- Row 10 : Test Case 6 --------------------------------------------------
$ws.Range("A10").Value = "Class:  gamble"
$ws.Range("C10").Value = "Joel: Test 7"
$ws.Range("E10").Value = "User types in ""!gamble tails"""
$ws.Range("F10").Value = "Bot tells the user whether they won the bet or not."
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = "Joel"
$ws.Range("J10").ClearContents()

# --- Row 11 : Test Case 7 ("Class: onSubHandler") ---------------------------
$ws.Range("A11").Value = "Class:  onSubHandler"
$ws.Range("C11").Value = "Joel: Test 8"
$ws.Range("E11").Value = "Wait for a subscription"
$ws.Range("F11").Value = "Message will appear in the command window saying ""/*/*/*/*/*Subscriber has been detected/*/*/*/*/*"""
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = "Joel"
$ws.Range("J11").ClearContents()

# --- Row 12 : Test Case 8 ---------------------------------------------------
$ws.Range("A12").Value = "Class:  onSubHandler"
$ws.Range("C12").Value = "Joel: Test 9"
$ws.Range("E12").Value = "User types ""/*/*/*/*/*Subscriber has been detected/*/*/*/*/*"""
$ws.Range("F12").Value = "No new sub will be detected"
$ws.Range("G12").Value = "Failed"
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = "Joel"
$ws.Range("J12").ClearContents()

# --- Cosmetic tweaks: column A widened, rows 4 & 12 taller for wrapped text
$ws.Columns("A").ColumnWidth = 23.67
$ws.Rows(4).RowHeight = 28.5
$ws.Rows(12).RowHeight = 26.25

# --- Leave the selection where the author last left it ---------------------
$ws.Range("F13").Select()
